# Load Screen code Updated
# - Replace the placeholder login e-mail in the "Login" sheet (B2) with the
#   new test address and turn it into a mailto hyperlink (Excel will add the
#   built-in "Hyperlink" style/font automatically).
# - Make "Login" the active sheet/tab again (it was "ResetPassword"),
#   leaving the selection on B11.

$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Login")

$wsLogin.Hyperlinks.Add($wsLogin.Range("B2"), "mailto:changepondtest10@yahoo.com", "", "", "changepondtest10@yahoo.com")

$wsLogin.Activate() | Out-Null
$wsLogin.Range("B11").Select() | Out-Null
